$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Regenerate orders: distances D51/D64/D80 -> D55/D69/D86 and size code S30 -> S31.
# Use whole-string replace semantics are not needed here because these tokens
# appear as substrings inside composed labels/filenames (e.g. "Face11_D51_S20",
# "Face11_D51_S20_l.png") as well as standalone ("D51", "S30"), so xlPart (2)
# LookAt is used to catch every occurrence, matching the shared-string diff.
$used.Replace("D51", "D55", 2, $null, $true) | Out-Null
$used.Replace("D64", "D69", 2, $null, $true) | Out-Null
$used.Replace("D80", "D86", 2, $null, $true) | Out-Null
$used.Replace("S30", "S31", 2, $null, $true) | Out-Null
